$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = [double]"1.77368233250472"
$ws.Range("B1").Value = [double]"0.002539946714806698"
$ws.Range("C1").Value = [double]"1.379315696670069"
$ws.Range("D1").Value = [double]"0.1889405185927815"
$ws.Range("E1").Value = [double]"1.570796395407675"
$ws.Range("F1").Value = [double]"0.2028859211499618"
$ws.Range("A2").Value = [double]"1.819293683204713"
$ws.Range("B2").Value = [double]"0.002378253947922411"
$ws.Range("C2").Value = [double]"1.380517981017872"
$ws.Range("D2").Value = [double]"0.1878999055330415"
$ws.Range("E2").Value = [double]"1.570796393670107"
$ws.Range("F2").Value = [double]"0.2484972728326442"
$ws.Range("A3").Value = [double]"2.023662773910207"
$ws.Range("B3").Value = [double]"0.001653763152656354"
$ws.Range("C3").Value = [double]"1.385905011990622"
$ws.Range("D3").Value = [double]"0.1832372691161953"
$ws.Range("E3").Value = [double]"1.570796385884652"
$ws.Range("F3").Value = [double]"0.4528663679412378"
$ws.Range("A4").Value = [double]"2.31119068080906"
$ws.Range("B4").Value = [double]"0.000634473390142702"
$ws.Range("C4").Value = [double]"1.39348405310025"
$ws.Range("D4").Value = [double]"0.1766773823694363"
$ws.Range("E4").Value = [double]"1.570796374931256"
$ws.Range("F4").Value = [double]"0.7403942810348323"
$ws.Range("A5").Value = [double]"2.515559771514555"
$ws.Range("B5").Value = [double]"-9.001740512335757e-05"
$ws.Range("C5").Value = [double]"1.398871084073"
$ws.Range("D5").Value = [double]"0.1720147459525901"
$ws.Range("E5").Value = [double]"1.570796367145801"
$ws.Range("F5").Value = [double]"0.9447633761434262"
$ws.Range("A6").Value = [double]"2.561171122214548"
$ws.Range("B6").Value = [double]"-0.000251710172007643"
$ws.Range("C6").Value = [double]"1.400073368420803"
$ws.Range("D6").Value = [double]"0.1709741328928501"
$ws.Range("E6").Value = [double]"1.570796365408233"
$ws.Range("F6").Value = [double]"0.9903747278261079"
